# Fruta / hortaliza, semanal
# Insert a new week's worth of price data (4 rows) for
# Vega Monumental Concepción - Pera, right after the existing row 188
# (before the old row 189), shifting all the following rows down by 4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 4 blank rows starting at row 189 (pushes old rows 189-250 to 193-254)
$ws.Range("A189:A192").EntireRow.Insert()

# --- New row 189: Packham's Triumph, Primera ---
$ws.Range("A189").Value2 = 11
$ws.Range("B189").Value2 = "Vega Monumental Concepción"
$ws.Range("C189").Value2 = "Bíobío"
$ws.Range("D189").Value2 = 44460
$ws.Range("E189").Value2 = 8
$ws.Range("F189").Value2 = "Fruta"
$ws.Range("G189").Value2 = 100104
$ws.Range("H189").Value2 = "Frutos de pepita"
$ws.Range("I189").Value2 = 100104005
$ws.Range("J189").Value2 = "Pera"
$ws.Range("K189").Value2 = "Packham's Triumph"
$ws.Range("L189").Value2 = "Primera"
$ws.Range("M189").Value2 = 200
$ws.Range("N189").Value2 = 10000
$ws.Range("O189").Value2 = 11000
$ws.Range("P189").Value2 = 10500
$ws.Range("Q189").Value2 = "`$/caja 16 kilos empedrada"
$ws.Range("R189").Value2 = "Región de O'Higgins"
$ws.Range("S189").Value2 = 656
$ws.Range("T189").Value2 = 16

# --- New row 190: Packham's Triumph, Segunda ---
$ws.Range("A190").Value2 = 11
$ws.Range("B190").Value2 = "Vega Monumental Concepción"
$ws.Range("C190").Value2 = "Bíobío"
$ws.Range("D190").Value2 = 44460
$ws.Range("E190").Value2 = 8
$ws.Range("F190").Value2 = "Fruta"
$ws.Range("G190").Value2 = 100104
$ws.Range("H190").Value2 = "Frutos de pepita"
$ws.Range("I190").Value2 = 100104005
$ws.Range("J190").Value2 = "Pera"
$ws.Range("K190").Value2 = "Packham's Triumph"
$ws.Range("L190").Value2 = "Segunda"
$ws.Range("M190").Value2 = 100
$ws.Range("N190").Value2 = 9000
$ws.Range("O190").Value2 = 9000
$ws.Range("P190").Value2 = 9000
$ws.Range("Q190").Value2 = "`$/caja 16 kilos empedrada"
$ws.Range("R190").Value2 = "Región de O'Higgins"
$ws.Range("S190").Value2 = 562
$ws.Range("T190").Value2 = 16

# --- New row 191: Winter Nelis, Primera ---
$ws.Range("A191").Value2 = 11
$ws.Range("B191").Value2 = "Vega Monumental Concepción"
$ws.Range("C191").Value2 = "Bíobío"
$ws.Range("D191").Value2 = 44460
$ws.Range("E191").Value2 = 8
$ws.Range("F191").Value2 = "Fruta"
$ws.Range("G191").Value2 = 100104
$ws.Range("H191").Value2 = "Frutos de pepita"
$ws.Range("I191").Value2 = 100104005
$ws.Range("J191").Value2 = "Pera"
$ws.Range("K191").Value2 = "Winter Nelis"
$ws.Range("L191").Value2 = "Primera"
$ws.Range("M191").Value2 = 50
$ws.Range("N191").Value2 = 10000
$ws.Range("O191").Value2 = 10000
$ws.Range("P191").Value2 = 10000
$ws.Range("Q191").Value2 = "`$/caja 16 kilos empedrada"
$ws.Range("R191").Value2 = "Región de O'Higgins"
$ws.Range("S191").Value2 = 625
$ws.Range("T191").Value2 = 16

# --- New row 192: Winter Nelis, Segunda ---
$ws.Range("A192").Value2 = 11
$ws.Range("B192").Value2 = "Vega Monumental Concepción"
$ws.Range("C192").Value2 = "Bíobío"
$ws.Range("D192").Value2 = 44460
$ws.Range("E192").Value2 = 8
$ws.Range("F192").Value2 = "Fruta"
$ws.Range("G192").Value2 = 100104
$ws.Range("H192").Value2 = "Frutos de pepita"
$ws.Range("I192").Value2 = 100104005
$ws.Range("J192").Value2 = "Pera"
$ws.Range("K192").Value2 = "Winter Nelis"
$ws.Range("L192").Value2 = "Segunda"
$ws.Range("M192").Value2 = 50
$ws.Range("N192").Value2 = 9000
$ws.Range("O192").Value2 = 9000
$ws.Range("P192").Value2 = 9000
$ws.Range("Q192").Value2 = "`$/caja 16 kilos empedrada"
$ws.Range("R192").Value2 = "Región de O'Higgins"
$ws.Range("S192").Value2 = 562
$ws.Range("T192").Value2 = 16
